$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 date value (precision correction)
$ws.Range("A2").Value = 45804.48548513889

# Insert new row 3 with updated price data
$ws.Range("A3").Value = 45804.43437102554
$ws.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B3").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C3").Value = "1Kg"
$ws.Range("D3").Value = "12,88€"
